# Applies the Yojimbo_Profits.xlsx value updates described in the commit diff.
# Each row below corresponds to one <row> element whose H:N (or subset) numeric
# cells were refreshed by the scheduled runner; some rows also gained or lost a
# trailing LeveProfitNQ/LeveProfitHQ cell (M/N) because the profit became (non-)negative
# in only one of the NQ/HQ variants.

$wb = $excel.ActiveWorkbook

# --- ALC sheet (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)

# Row 19
$ws.Range("H19").Value = 800.5
$ws.Range("I19").Value = 800.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 800.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -625.5
$ws.Range("N19").ClearContents()

# Row 117
$ws.Range("H117").Value = 40618.332
$ws.Range("J117").Value = 40618.332
$ws.Range("L117").Value = 40618.332
$ws.Range("N117").Value = -49796.332

# Row 128
$ws.Range("H128").Value = 43621.668
$ws.Range("J128").Value = 43621.668
$ws.Range("L128").Value = 43621.668
$ws.Range("N128").Value = -53581.668

# Row 133
$ws.Range("H133").Value = 76593.336
$ws.Range("J133").Value = 76593.336
$ws.Range("L133").Value = 76593.336
$ws.Range("N133").Value = -86713.336

# Row 137
$ws.Range("H137").Value = 6317.8237
$ws.Range("I137").Value = 6683.3335
$ws.Range("J137").Value = 5440.6
$ws.Range("K137").Value = 20050.0005
$ws.Range("L137").Value = 16321.8
$ws.Range("M137").Value = -17500.0005
$ws.Range("N137").Value = -21421.8

# Row 140
$ws.Range("H140").Value = 40940
$ws.Range("J140").Value = 40940
$ws.Range("L140").Value = 40940
$ws.Range("N140").Value = -51300


# --- ARM sheet (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)

# Row 61
$ws.Range("H61").Value = 5584.1816
$ws.Range("I61").Value = 5976.9
$ws.Range("J61").Value = 1657
$ws.Range("K61").Value = 5976.9
$ws.Range("L61").Value = 1657
$ws.Range("M61").Value = -5764.9
$ws.Range("N61").Value = -2081

# Row 74
$ws.Range("H74").Value = 4186.3335
$ws.Range("I74").Value = 4428.222
$ws.Range("J74").Value = 2009.3334
$ws.Range("K74").Value = 4428.222
$ws.Range("L74").Value = 2009.3334
$ws.Range("M74").Value = -3554.222
$ws.Range("N74").Value = -3757.3334

# Row 77
$ws.Range("H77").Value = 4186.3335
$ws.Range("I77").Value = 4428.222
$ws.Range("J77").Value = 2009.3334
$ws.Range("K77").Value = 22141.11
$ws.Range("L77").Value = 10046.667
$ws.Range("M77").Value = -17773.11
$ws.Range("N77").Value = -18782.667

# Row 122
$ws.Range("H122").Value = 2005.8
$ws.Range("I122").Value = 2077.7144
$ws.Range("J122").Value = 1838
$ws.Range("K122").Value = 6233.1432
$ws.Range("L122").Value = 5514
$ws.Range("M122").Value = -3783.1432
$ws.Range("N122").Value = -10414

# Row 123
$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800

# Row 132
$ws.Range("H132").Value = 13765.529
$ws.Range("I132").Value = 9308.923000000001
$ws.Range("J132").Value = 28249.5
$ws.Range("K132").Value = 27926.769
$ws.Range("L132").Value = 84748.5
$ws.Range("M132").Value = -25396.769
$ws.Range("N132").Value = -89808.5

# Row 136
$ws.Range("H136").Value = 5584.1816
$ws.Range("I136").Value = 5976.9
$ws.Range("J136").Value = 1657
$ws.Range("K136").Value = 17930.7
$ws.Range("L136").Value = 4971
$ws.Range("M136").Value = -15380.7
$ws.Range("N136").Value = -10071


# --- BSM sheet (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)

# Row 64
$ws.Range("H64").Value = 815.4
$ws.Range("J64").Value = 856.2857
$ws.Range("L64").Value = 856.2857
$ws.Range("N64").Value = -1306.2857

# Row 67
$ws.Range("H67").Value = 815.4
$ws.Range("J67").Value = 856.2857
$ws.Range("L67").Value = 856.2857
$ws.Range("N67").Value = -2416.2857

# Row 99
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3502
$ws.Range("N99").ClearContents()

# Row 134
$ws.Range("H134").Value = 3096.5957
$ws.Range("I134").Value = 3105.3333
$ws.Range("J134").Value = 2900
$ws.Range("K134").Value = 9315.999899999999
$ws.Range("L134").Value = 8700
$ws.Range("M134").Value = -6780.999899999999
$ws.Range("N134").Value = -13770


# --- CRP sheet (Worksheets.Item(4)) ---
$ws = $wb.Worksheets.Item(4)

# Row 31
$ws.Range("H31").Value = 22853.559
$ws.Range("I31").Value = 32498.943
$ws.Range("J31").Value = 2995.4119
$ws.Range("K31").Value = 32498.943
$ws.Range("L31").Value = 2995.4119
$ws.Range("M31").Value = -32203.943
$ws.Range("N31").Value = -3585.4119

# Row 34
$ws.Range("H34").Value = 22853.559
$ws.Range("I34").Value = 32498.943
$ws.Range("J34").Value = 2995.4119
$ws.Range("K34").Value = 32498.943
$ws.Range("L34").Value = 2995.4119
$ws.Range("M34").Value = -32296.943
$ws.Range("N34").Value = -3399.4119

# Row 58
$ws.Range("H58").Value = 1368.0714
$ws.Range("I58").Value = 1518.5834
$ws.Range("J58").Value = 465
$ws.Range("K58").Value = 1518.5834
$ws.Range("L58").Value = 465
$ws.Range("M58").Value = -1315.5834
$ws.Range("N58").Value = -871

# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Row 132
$ws.Range("H132").Value = 8222.666999999999
$ws.Range("I132").Value = 9000.786
$ws.Range("J132").Value = 5499.25
$ws.Range("K132").Value = 27002.358
$ws.Range("L132").Value = 16497.75
$ws.Range("M132").Value = -24472.358
$ws.Range("N132").Value = -21557.75

# Row 134
$ws.Range("H134").Value = 8744.464
$ws.Range("I134").Value = 6505.5713
$ws.Range("J134").Value = 15461.143
$ws.Range("K134").Value = 19516.7139
$ws.Range("L134").Value = 46383.429
$ws.Range("M134").Value = -16981.7139
$ws.Range("N134").Value = -51453.429

# Row 136
$ws.Range("H136").Value = 1368.0714
$ws.Range("I136").Value = 1518.5834
$ws.Range("J136").Value = 465
$ws.Range("K136").Value = 4555.7502
$ws.Range("L136").Value = 1395
$ws.Range("M136").Value = -2005.7502
$ws.Range("N136").Value = -6495


# --- GSM sheet (Worksheets.Item(6)) ---
$ws = $wb.Worksheets.Item(6)

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 102
$ws.Range("H102").Value = 3500
$ws.Range("I102").Value = 3500
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3500
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1878
$ws.Range("N102").ClearContents()

# Row 122
$ws.Range("H122").Value = 1272.0952
$ws.Range("I122").Value = 1272.0952
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3816.2856
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1366.2856
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 10323.044
$ws.Range("I132").Value = 6421.55
$ws.Range("J132").Value = 36333
$ws.Range("K132").Value = 19264.65
$ws.Range("L132").Value = 108999
$ws.Range("M132").Value = -16734.65
$ws.Range("N132").Value = -114059


# --- LTW sheet (Worksheets.Item(7)) ---
$ws = $wb.Worksheets.Item(7)

# Row 122
$ws.Range("H122").Value = 4334.1665
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4501.25
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 13503.75
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -18403.75

# Row 132
$ws.Range("H132").Value = 44054.418
$ws.Range("I132").Value = 68260.336
$ws.Range("J132").Value = 3711.2222
$ws.Range("K132").Value = 204781.008
$ws.Range("L132").Value = 11133.6666
$ws.Range("M132").Value = -202251.008
$ws.Range("N132").Value = -16193.6666

# Row 136
$ws.Range("H136").Value = 2393.7659
$ws.Range("I136").Value = 1884.7333
$ws.Range("J136").Value = 3292.0588
$ws.Range("K136").Value = 5654.199900000001
$ws.Range("L136").Value = 9876.1764
$ws.Range("M136").Value = -3104.2001
$ws.Range("N136").Value = -14976.1764


# --- WVR sheet (Worksheets.Item(8)) ---
$ws = $wb.Worksheets.Item(8)

# Row 92
$ws.Range("H92").Value = 34367
$ws.Range("J92").Value = 34367
$ws.Range("L92").Value = 34367
$ws.Range("N92").Value = -39359

# Row 122
$ws.Range("H122").Value = 11177713
$ws.Range("I122").Value = 12692307
$ws.Range("K122").Value = 38076921
$ws.Range("M122").Value = -38074471

# Row 123
$ws.Range("H123").Value = 19500
$ws.Range("J123").Value = 19500
$ws.Range("L123").Value = 19500
$ws.Range("N123").Value = -29300

# Row 132
$ws.Range("H132").Value = 4458.75
$ws.Range("I132").Value = 4589.8125
$ws.Range("J132").Value = 4109.25
$ws.Range("K132").Value = 13769.4375
$ws.Range("L132").Value = 12327.75
$ws.Range("M132").Value = -11239.4375
$ws.Range("N132").Value = -17387.75

# Row 136
$ws.Range("H136").Value = 24393960
$ws.Range("I136").Value = 31254090
$ws.Range("J136").Value = 2388.3333
$ws.Range("K136").Value = 93762270
$ws.Range("L136").Value = 7164.999899999999
$ws.Range("M136").Value = -93759720
$ws.Range("N136").Value = -12264.9999

